# Applies the "Astronomy" -> "Biology" rewrite described in the commit diff.
#
# Word's Find.Execute(..., Replace:=wdReplaceOne/All) path runs the replacement
# text through AutoCorrect (straight "'" -> curly "'"), which the source diff
# does not want. To keep a literal apostrophe we instead use Find.Execute only
# to *locate* the target span (Replace:=wdReplaceNone) and then assign
# Range.Text directly, which performs a verbatim, format-preserving swap.

function Replace-FirstMatch {
    param(
        [string]$OldText,
        [string]$NewText
    )
    $rng = $d.Content
    $found = $rng.Find.Execute($OldText, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        throw "Could not find text: $OldText"
    }
    $rng.Text = $NewText
}

$d = $word.ActiveDocument

# --- Title ---------------------------------------------------------------
Replace-FirstMatch "Astronomy's Eye on Celestial Symphony" `
    "The Essence of Biology Unveiled: Exploring the Realm of Life's Complexities"

# --- Byline ----------------------------------------------------------------
Replace-FirstMatch "Isabella Matthews" "Amelia Williams"

# --- Email (collapses the isabella / . / matthews@spatiumobservatory runs) -
Replace-FirstMatch "isabella.matthews@spatiumobservatory" "amwilliams97@educonnect"

# --- Body paragraph ----------------------------------------------------------
Replace-FirstMatch `
    "Across the vast panorama of the cosmos, where the reverie of stars enchants the minds of humankind, an intricate symphony of cosmic phenomena awaits our exploration" `
    "Biology, the study of life in all its intricate forms, offers an enthralling journey into the remarkable phenomena that surround us"

Replace-FirstMatch `
    " Immerse yourself in this celestial ballet, where galaxies pirouette, planets waltz in orbital harmony, and cosmic winds hum through interstellar voids" `
    " From the swirling patterns of a single cell to the sophisticated symphony of ecosystems, biology unveils the secrets of life's boundless diversity and its inherent unity"

Replace-FirstMatch `
    " From the titanic dance of black holes to the radiant birth and death of stars, the universe hums with untold stories. Humans have always looked towards the sky with wonder, curiosity, and a yearning to comprehend the cosmos. In this symphony of celestial wonders, we find a tapestry of knowledge that shapes our understanding of existence, time, and our place in the grand symphony of the universe" `
    " In this realm, we delve into the captivating secrets of living organisms - their structures, intricate processes, and interactions with one another and the wider world"

Replace-FirstMatch `
    "Our journey through this celestial symphony begins with the majestic dance of galaxies, spiraled wonders whirling in their gravitational embrace" `
    "Biology is a captivating mosaic of interconnected concepts, a captivating saga of life's astounding resilience"

Replace-FirstMatch `
    " Like graceful ballerinas adoring the cosmic stage, they waltz across vast distances, their intricate choreography revealing the mysteries of gravity and dark matter" `
    " Witnessing the transformation of a caterpillar into a butterfly or unraveling the complex chain of events that allows plants to turn sunlight into energy illuminates the sheer mindboggling capabilities of life"

Replace-FirstMatch `
    " As we delve deeper into this celestial waltz, we encounter the captivating ballet of planetary motion, revolving around their stellar partners in an eternal dance of gravitational allure. Intriguing exoplanet discoveries paint a portrait of diverse worlds, ranging from colossal gas giants to terrestrial planets, each harboring secrets waiting to be unraveled" `
    " Biology incites in us an innate quest for answers, unraveling the mysteries of our bodies, understanding the mechanisms behind diseases, and searching for innovative solutions to address pressing ecological issues"

Replace-FirstMatch `
    "Further, this celestial symphony orchestrates a mesmerizing array of celestial spectacles. Meteors streak across the black velvet curtain of the night sky, leaving ephemeral trails of luminescent beauty. Exploding stars erupt in radiant finales, illuminating entire galaxies in their fiery brilliance. Supernovas, akin to cosmic fireworks, forge the elements that shape the universe, while black holes lurk as enigmatic conductors, warbling distorted notes of gravity, and devouring matter, adding a haunting beauty to the cosmic symphony" `
    "As we unravel the intricate tapestry of life, we come face-to-face with questions that have pondered humanity for eons: How did life originate? How do organisms adapt and evolve in response to their surroundings? How can we decipher the genetic language that dictates the symphony of life? Biology provides us with a framework to tackle these perplexing inquiries, inviting us to unlock the secrets of our existence"

# --- Summary paragraph -------------------------------------------------------
Replace-FirstMatch `
    "In the interwoven tapestry of cosmic phenomena, humanity finds a symphony of wonder, knowledge, and boundless exploration" `
    "The study of biology unveils the captivating intricacies of life, from the minuscule world within a single cell to the complex interactions of entire ecosystems"

Replace-FirstMatch `
    " This celestial ballet invites introspection, inspiring us to contemplate our place in the universe and our connection to the grand orchestration of existence" `
    " It encapsulates investigations into diverse living organisms, encompassing their structures, functions, and interplay with each other and the environment"

Replace-FirstMatch `
    " Astronomy's eye unveils the intricate dance of galaxies, the waltz of planets, and the majestic spectacle of celestial events, painting a portrait of the cosmos that is both awe-inspiring and profound" `
    " Biology inspires us to delve into profound questions about the origins of life, mechanisms of adaptation, and the intricate genetic code"

Replace-FirstMatch `
    " Our exploration of this celestial symphony is an odyssey through the vastness of space and time, offering perspectives that transcend earthly boundaries and ignite imaginations" `
    " The journey of understanding biology unveils the essence of our connection with all living things and unravels the fascinating story of our place within the grand tapestry of life on Earth"

# --- Trailing empty paragraph added at the end of the document -------------
$d.Content.InsertParagraphAfter()
